# Generate Report for Handback
# - Mark every row's Status as "Handed back: in sync with en-US" (was "Ready for handoff")
# - Fill in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   columns on the per-locale sheets (zh-cn, de-de) now that the handback files exist
# - Widen a few columns that now hold longer content

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFile1 = "4feb7aad-90be-416b-9afe-7194f5e88201.md"
$mdFile2 = "b003a3cc-756c-4406-ab57-bbec474a61c0.md"
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d6c761c6b414eb92bbb6004edaa7acc73e8c890/e2e/4feb7aad-90be-416b-9afe-7194f5e88201.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d6c761c6b414eb92bbb6004edaa7acc73e8c890/e2e/b003a3cc-756c-4406-ab57-bbec474a61c0.md"

# ---------------------------------------------------------------------------
# Overview sheet: Status column (E, F) -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Range("I2").Value = $mdFile1
$wsZhCn.Range("J2").Value = "4feb7aad-90be-416b-9afe-7194f5e88201.1bc369c1e67c487b37f3400ca5226a774266676f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-01 01:05:54"

$wsZhCn.Range("I3").Value = $mdFile2
$wsZhCn.Range("J3").Value = "b003a3cc-756c-4406-ab57-bbec474a61c0.84b3f4a4a727a2e5672184347ddb206731031e31.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-01 01:05:54"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Range("I2").Value = $mdFile1
$wsDeDe.Range("J2").Value = "4feb7aad-90be-416b-9afe-7194f5e88201.1bc369c1e67c487b37f3400ca5226a774266676f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-01 01:06:06"

$wsDeDe.Range("I3").Value = $mdFile2
$wsDeDe.Range("J3").Value = "b003a3cc-756c-4406-ab57-bbec474a61c0.84b3f4a4a727a2e5672184347ddb206731031e31.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-01 01:06:06"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
